$wb = $excel.ActiveWorkbook

# --- Update "Last Updated" timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 04:05 PM"

# --- Update the "1 Month Performance" table (rows 68-74 shift up, new entry added) ---
$perf = $wb.Worksheets.Item("1 Month Performance")

$perf.Range("B68").Value = "RBLBANK"
$perf.Range("C68").Value = 19.2556

$perf.Range("B69").Value = "MOLDTECH"
$perf.Range("C69").Value = 19.1891

$perf.Range("B70").Value = "THOMASCOTT"
$perf.Range("C70").Value = 19.1649

$perf.Range("B71").Value = "KARURVYSYA"
$perf.Range("C71").Value = 19.11

$perf.Range("B72").Value = "IIFL"
$perf.Range("C72").Value = 18.9853

$perf.Range("B73").Value = "LUMAXIND"
$perf.Range("C73").Value = 18.8057

$perf.Range("B74").Value = "REPRO"
$perf.Range("C74").Value = 18.689
